$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row => (B, C, D)
$data = @{
  2  = @(100, 0.7785501480102539, 0)
  3  = @(100, 0.4487638473510742, 0)
  4  = @(0,   0.4951229095458984, 0)
  5  = @(0,   0.4708287715911865, 0)
  6  = @(100, 0.07204771041870117, 0)
  7  = @(100, 0.06608104705810547, 0)
  8  = @(0,   0.09218668937683105, 0)
  9  = @(0,   0.1927340030670166, 0)
  10 = @(100, 0.6884679794311523, 0)
  11 = @(100, 0.257838249206543, 0)
  12 = @(100, 0.5248048305511475, 0)
  13 = @(100, 0.3641648292541504, 0)
  14 = @(100, 0.4229428768157959, 0)
  15 = @(100, 0.4173502922058105, 0)
  16 = @(0,   0.5478470325469971, 0)
  17 = @(0,   0.5007138252258301, 0)
  18 = @(100, 0.05916213989257812, 0)
  19 = @(100, 0.0615692138671875, 0)
  20 = @(0,   0.1086909770965576, 0)
  21 = @(0,   0.09883594512939453, 0)
  22 = @(100, 0.2533378601074219, 0)
  23 = @(100, 0.2440938949584961, 0)
  24 = @(100, 0.3615798950195312, 0)
  25 = @(100, 0.4897820949554443, 0)
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Cells.Item($row, 2).Value = $vals[0]
  $ws.Cells.Item($row, 3).Value = $vals[1]
  $ws.Cells.Item($row, 4).Value = $vals[2]
}
